$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Unmerge A1:F2 and re-merge as A1:G2 so the new column G is included in
#    the title banner merge (matches mergeCell order/ref in target).
# ---------------------------------------------------------------------------
$ws.Range("A1:F2").UnMerge()

# ---------------------------------------------------------------------------
# 2) Give row 1 and row 2 an explicit custom height (15pt) like the target.
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 15
$ws.Rows(2).RowHeight = 15

$ws.Range("A1:G2").Merge()

# ---------------------------------------------------------------------------
# 3) Add the new "Revisado" header cell in G4, matching the style used for
#    the other header cells (left aligned) plus a new light-green fill.
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = "Revisado"
$ws.Range("G4").Interior.Color = 9359785

# ---------------------------------------------------------------------------
# 4) Populate the new "Revisado" date column (G) for the data rows, copying
#    the date format used by column F, and leaving some rows blank.
# ---------------------------------------------------------------------------
$datedRows = @(6,7,9,11,12,13,14,15,18,21,26,27)
foreach ($r in $datedRows) {
    $ws.Range("F$r").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)
    $ws.Range("G$r").Value2 = 45759
}

$ws.Range("F19").Copy()
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("G19").Value2 = 45766

$blankRows = @(8,10,16,17,20,22,23,24,25)
foreach ($r in $blankRows) {
    $ws.Range("F$r").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 5) Make column G the same (best-fit-like) width as column F.
# ---------------------------------------------------------------------------
$ws.Columns("G:G").AutoFit()

# ---------------------------------------------------------------------------
# 6) Update the sheet view: scroll back to the top, select I24 (matches the
#    saved selection/topLeftCell in the target file).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("I24").Select()

# ---------------------------------------------------------------------------
# 7) Page setup - paper size / orientation explicitly recorded on save.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
